# Apply the benchmark-stat corrections described by the commit:
# "Fixed README.md stats and docx preparation for all Renaissance -
#  JDK 17 - Shenandoah GC tests"
#
# The document is a single-column, 46-row table. The first dozen rows
# hold per-metric summary values, and the last three rows hold
# multi-column (tab-separated) raw data that collapses down to a
# single summary figure.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $newText) {
    $cell = $table.Cell($row, 1)
    $cell.Range.Text = $newText
}

# Simple value swaps in the summary rows near the top of the table.
Set-CellText $t 1  "0M"
Set-CellText $t 2  "0M"
Set-CellText $t 3  "0M"
Set-CellText $t 4  "1514"
Set-CellText $t 5  "0.00001"
Set-CellText $t 6  "0.00085"
Set-CellText $t 7  "0.00011"
Set-CellText $t 9  "0.00014"
Set-CellText $t 10 "0.00015"
Set-CellText $t 11 "0.00017"
Set-CellText $t 12 "0.17944"

# The three raw, tab-separated detail rows near the bottom collapse
# down to a single value each.
Set-CellText $t 44 "99.91"
Set-CellText $t 45 "0.18"
Set-CellText $t 46 "203"
